$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Program names (column B) to the new, more descriptive titles.
$ws.Range("B2").Value  = "National Water Resource Management Sector Framework"
$ws.Range("B3").Value  = "Overarching Flood Risk Management Framework"
$ws.Range("B4").Value  = "Overarching National Drought Risk Management Framework"
$ws.Range("B18").Value = "National Disaster Risk Management Sector Framework"
$ws.Range("B19").Value = "Overarching Flood Risk Management Framework"
$ws.Range("B32").Value = "Overarching National Drought Risk Management Framework"

# Update the active selection / scroll position to match the saved view.
$ws.Range("B18").Select() | Out-Null
